$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44530
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 200

# Row 3
$ws.Cells.Item(3, 4).Value = 44530
$ws.Cells.Item(3, 12).Value = 'Segunda'
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 16000
$ws.Cells.Item(3, 16).Value = 16000
$ws.Cells.Item(3, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(3, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(3, 19).Value = 2000
$ws.Cells.Item(3, 20).Value = 8

# Row 4
$ws.Cells.Item(4, 4).Value = 44512
$ws.Cells.Item(4, 12).Value = 'Segunda'
$ws.Cells.Item(4, 13).Value = 300
$ws.Cells.Item(4, 14).Value = 19000
$ws.Cells.Item(4, 15).Value = 20000
$ws.Cells.Item(4, 16).Value = 19500
$ws.Cells.Item(4, 19).Value = 2438

# Row 5
$ws.Cells.Item(5, 4).Value = 45247
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 14).Value = 25000
$ws.Cells.Item(5, 15).Value = 25000
$ws.Cells.Item(5, 16).Value = 25000
$ws.Cells.Item(5, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(5, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(5, 19).Value = 2500
$ws.Cells.Item(5, 20).Value = 10

# Row 6
$ws.Cells.Item(6, 4).Value = 44533
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 300

# Row 7
$ws.Cells.Item(7, 4).Value = 44533
$ws.Cells.Item(7, 12).Value = 'Segunda'
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 16000
$ws.Cells.Item(7, 15).Value = 16000
$ws.Cells.Item(7, 16).Value = 16000
$ws.Cells.Item(7, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(7, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(7, 19).Value = 2000
$ws.Cells.Item(7, 20).Value = 8

# Row 8
$ws.Cells.Item(8, 4).Value = 44523
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 400
$ws.Cells.Item(8, 14).Value = 21000
$ws.Cells.Item(8, 15).Value = 22000
$ws.Cells.Item(8, 16).Value = 21500
$ws.Cells.Item(8, 19).Value = 2688

# Row 9
$ws.Cells.Item(9, 4).Value = 44523
$ws.Cells.Item(9, 13).Value = 100
$ws.Cells.Item(9, 14).Value = 18000
$ws.Cells.Item(9, 15).Value = 18000
$ws.Cells.Item(9, 16).Value = 18000
$ws.Cells.Item(9, 19).Value = 2250

# Row 10
$ws.Cells.Item(10, 4).Value = 44162
$ws.Cells.Item(10, 14).Value = 2000
$ws.Cells.Item(10, 15).Value = 2100
$ws.Cells.Item(10, 16).Value = 2050
$ws.Cells.Item(10, 17).Value = '$/kilo (en caja de 14 kilos)'
$ws.Cells.Item(10, 19).Value = 2050
$ws.Cells.Item(10, 20).Value = 1

# Row 11
$ws.Cells.Item(11, 4).Value = 45240
$ws.Cells.Item(11, 13).Value = 150
$ws.Cells.Item(11, 14).Value = 28000
$ws.Cells.Item(11, 15).Value = 28000
$ws.Cells.Item(11, 16).Value = 28000
$ws.Cells.Item(11, 19).Value = 2800

# Row 12
$ws.Cells.Item(12, 4).Value = 44519
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 400
$ws.Cells.Item(12, 14).Value = 21000
$ws.Cells.Item(12, 15).Value = 22000
$ws.Cells.Item(12, 16).Value = 21500
$ws.Cells.Item(12, 19).Value = 2688

# Row 13
$ws.Cells.Item(13, 4).Value = 44519
$ws.Cells.Item(13, 12).Value = 'Segunda'
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 18000
$ws.Cells.Item(13, 15).Value = 18000
$ws.Cells.Item(13, 16).Value = 18000
$ws.Cells.Item(13, 19).Value = 2250

# Row 14
$ws.Cells.Item(14, 4).Value = 44498
$ws.Cells.Item(14, 13).Value = 300
$ws.Cells.Item(14, 14).Value = 19000
$ws.Cells.Item(14, 15).Value = 20000
$ws.Cells.Item(14, 16).Value = 19500
$ws.Cells.Item(14, 19).Value = 2438

# Row 15
$ws.Cells.Item(15, 4).Value = 44516
$ws.Cells.Item(15, 12).Value = 'Segunda'
$ws.Cells.Item(15, 14).Value = 18000
$ws.Cells.Item(15, 15).Value = 19000
$ws.Cells.Item(15, 16).Value = 18500
$ws.Cells.Item(15, 19).Value = 2312

# Row 16
$ws.Cells.Item(16, 4).Value = 44890
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 22000
$ws.Cells.Item(16, 15).Value = 22500
$ws.Cells.Item(16, 16).Value = 22250
$ws.Cells.Item(16, 19).Value = 2781

# Row 17
$ws.Cells.Item(17, 4).Value = 44495
$ws.Cells.Item(17, 13).Value = 270
$ws.Cells.Item(17, 14).Value = 19000
$ws.Cells.Item(17, 15).Value = 20000
$ws.Cells.Item(17, 16).Value = 19556
$ws.Cells.Item(17, 19).Value = 2444

# Row 18
$ws.Cells.Item(18, 4).Value = 44873
$ws.Cells.Item(18, 13).Value = 300

# Row 19
$ws.Cells.Item(19, 4).Value = 44505
$ws.Cells.Item(19, 12).Value = 'Segunda'
$ws.Cells.Item(19, 13).Value = 300
$ws.Cells.Item(19, 14).Value = 19000
$ws.Cells.Item(19, 15).Value = 20000
$ws.Cells.Item(19, 16).Value = 19500
$ws.Cells.Item(19, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(19, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(19, 19).Value = 2438
$ws.Cells.Item(19, 20).Value = 8

# Row 20
$ws.Cells.Item(20, 4).Value = 44159
$ws.Cells.Item(20, 14).Value = 2000
$ws.Cells.Item(20, 15).Value = 2100
$ws.Cells.Item(20, 16).Value = 2050
$ws.Cells.Item(20, 17).Value = '$/kilo (en caja de 14 kilos)'
$ws.Cells.Item(20, 19).Value = 2050
$ws.Cells.Item(20, 20).Value = 1

# Row 21
$ws.Cells.Item(21, 4).Value = 45237
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 150
$ws.Cells.Item(21, 14).Value = 28000
$ws.Cells.Item(21, 15).Value = 28000
$ws.Cells.Item(21, 16).Value = 28000
$ws.Cells.Item(21, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(21, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(21, 19).Value = 2800
$ws.Cells.Item(21, 20).Value = 10

# Row 22
$ws.Cells.Item(22, 4).Value = 44526
$ws.Cells.Item(22, 12).Value = 'Primera'
$ws.Cells.Item(22, 13).Value = 300
$ws.Cells.Item(22, 14).Value = 21000
$ws.Cells.Item(22, 15).Value = 21000
$ws.Cells.Item(22, 16).Value = 21000
$ws.Cells.Item(22, 19).Value = 2625

# Row 23
$ws.Cells.Item(23, 4).Value = 44488
$ws.Cells.Item(23, 12).Value = 'Segunda'
$ws.Cells.Item(23, 13).Value = 160
$ws.Cells.Item(23, 14).Value = 17000
$ws.Cells.Item(23, 15).Value = 18000
$ws.Cells.Item(23, 16).Value = 17500
$ws.Cells.Item(23, 19).Value = 2188

# Row 24
$ws.Cells.Item(24, 4).Value = 45244
$ws.Cells.Item(24, 13).Value = 150
$ws.Cells.Item(24, 14).Value = 25000
$ws.Cells.Item(24, 15).Value = 25000
$ws.Cells.Item(24, 16).Value = 25000
$ws.Cells.Item(24, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(24, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(24, 19).Value = 2500
$ws.Cells.Item(24, 20).Value = 10

# Row 25
$ws.Cells.Item(25, 4).Value = 44876

# Row 26
$ws.Cells.Item(26, 4).Value = 44491
$ws.Cells.Item(26, 12).Value = 'Segunda'
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 18000
$ws.Cells.Item(26, 15).Value = 19000
$ws.Cells.Item(26, 16).Value = 18500
$ws.Cells.Item(26, 19).Value = 2312

# Row 27
$ws.Cells.Item(27, 4).Value = 44880
$ws.Cells.Item(27, 13).Value = 300
$ws.Cells.Item(27, 14).Value = 22000
$ws.Cells.Item(27, 15).Value = 22500
$ws.Cells.Item(27, 16).Value = 22250
$ws.Cells.Item(27, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(27, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(27, 19).Value = 2781
$ws.Cells.Item(27, 20).Value = 8

# Row 28
$ws.Cells.Item(28, 4).Value = 44509
$ws.Cells.Item(28, 13).Value = 200
$ws.Cells.Item(28, 16).Value = 19500
$ws.Cells.Item(28, 19).Value = 2438

# Row 29
$ws.Cells.Item(29, 4).Value = 44894
$ws.Cells.Item(29, 14).Value = 22000
$ws.Cells.Item(29, 15).Value = 22500
$ws.Cells.Item(29, 16).Value = 22250
$ws.Cells.Item(29, 19).Value = 2781

# Row 30
$ws.Cells.Item(30, 4).Value = 44895
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 200
$ws.Cells.Item(30, 14).Value = 22000
$ws.Cells.Item(30, 15).Value = 22500
$ws.Cells.Item(30, 16).Value = 22250
$ws.Cells.Item(30, 19).Value = 2781

